$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '50.881.39'
$ws.Range('E2').Value = '  -2.63%  '
$ws.Range('D3').Value = '2.883.72'
$ws.Range('E3').Value = '  -2.73%  '
$ws.Range('E4').Value = '  -0.23%  '
$ws.Range('D5').Value = '''365.39'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +2.42%  '
$ws.Range('D6').Value = '''101.38'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -7.38%  '
$ws.Range('D7').Value = '''0.536'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -5.71%  '
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('D9').Value = '''0.582'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -7.04%  '
$ws.Range('D10').Value = '''36.31'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -6.72%  '
$ws.Range('E11').Value = '  +0.60%  '
$ws.Range('D12').Value = '''0.0828'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -4.73%  '
$ws.Range('D13').Value = '''18.21'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -5.97%  '
$ws.Range('D14').Value = '3.345.58'
$ws.Range('E14').Value = '  -3.37%  '
$ws.Range('D15').Value = '''7.33'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -6.16%  '
$ws.Range('D16').Value = '2.880.66'
$ws.Range('E16').Value = '  -5.17%  '
$ws.Range('D17').Value = '''0.926'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -5.95%  '
$ws.Range('D18').Value = '50.842.31'
$ws.Range('E18').Value = '  -2.54%  '
$ws.Range('D19').Value = '''3.24'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -7.61%  '
$ws.Range('D20').Value = '''7.15'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -6.22%  '
$ws.Range('D21').Value = '''12.71'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -8.07%  '
$ws.Range('D22').Value = '0.0₃0934'
$ws.Range('E22').Value = '  -4.85%  '
$ws.Range('D23').Value = '''67.76'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -3.70%  '
$ws.Range('D24').Value = '''256.85'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -4.49%  '
$ws.Range('D25').Value = '''2.66'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -4.79%  '
$ws.Range('B26').Value = 'Dai'
$ws.Range('C26').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D26').Value = '''1.00'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +0.06%  '
$ws.Range('B27').Value = 'Kaspa'
$ws.Range('C27').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D27').Value = '''0.167'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -5.88%  '
$ws.Range('B28').Value = 'EthereumClassic'
$ws.Range('C28').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D28').Value = '''25.46'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -6.29%  '
$ws.Range('B29').Value = 'Filecoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D29').Value = '''6.88'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -9.48%  '
$ws.Range('B30').Value = 'Hedera'
$ws.Range('C30').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D30').Value = '''0.102'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -4.81%  '
$ws.Range('B31').Value = 'Cosmos'
$ws.Range('C31').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D31').Value = '''9.79'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -6.13%  '
$ws.Range('B32').Value = 'RenderToken'
$ws.Range('C32').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D32').Value = '''5.92'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -2.59%  '
$ws.Range('B33').Value = 'Toncoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D33').Value = '''2.13'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -2.83%  '
$ws.Range('B34').Value = 'InjectiveProtocol'
$ws.Range('C34').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D34').Value = '''34.18'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -8.82%  '
$ws.Range('B35').Value = 'OKB'
$ws.Range('C35').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D35').Value = '''50.76'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -2.69%  '
$ws.Range('B36').Value = 'FirstDigitalUSD'
$ws.Range('C36').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D36').Value = '''1.00'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +0.25%  '
$ws.Range('B37').Value = 'VeChain'
$ws.Range('C37').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D37').Value = '''0.0413'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -6.63%  '
$ws.Range('B38').Value = 'LidoDAOToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D38').Value = '''3.01'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -6.63%  '
$ws.Range('B39').Value = 'Stacks'
$ws.Range('C39').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D39').Value = '''2.60'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -4.53%  '
$ws.Range('B40').Value = 'Celestia'
$ws.Range('C40').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D40').Value = '''16.79'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -7.13%  '
$ws.Range('B41').Value = 'ARBITRUM'
$ws.Range('C41').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D41').Value = '''1.83'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -10.24%  '
$ws.Range('B42').Value = 'Stellar'
$ws.Range('C42').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D42').Value = '''0.112'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -6.15%  '
$ws.Range('B43').Value = 'EnergySwap'
$ws.Range('C43').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D43').Value = '''21.76'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -7.61%  '
$ws.Range('B44').Value = 'Monero'
$ws.Range('C44').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D44').Value = '''117.49'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -1.43%  '
$ws.Range('B45').Value = 'WEMIXToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D45').Value = '''2.09'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -3.75%  '
$ws.Range('B46').Value = 'Maker'
$ws.Range('C46').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D46').Value = '2.018.99'
$ws.Range('E46').Value = '  -5.71%  '
$ws.Range('B47').Value = 'ApeXProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D47').Value = '''2.31'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -6.44%  '
$ws.Range('B48').Value = 'NEARProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D48').Value = '''3.11'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -9.82%  '
$ws.Range('B49').Value = 'RocketPoolETH'
$ws.Range('C49').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D49').Value = '3.185.78'
$ws.Range('E49').Value = '  -2.84%  '
$ws.Range('B50').Value = 'TheGraph'
$ws.Range('C50').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D50').Value = '''0.233'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -4.25%  '
$ws.Range('B51').Value = 'BEAM'
$ws.Range('C51').Value = 'https://coinranking.com/coin/cYYMfXF4u+beam-beam'
$ws.Range('D51').Value = '''0.0308'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -12.34%  '
